$wb = $excel.ActiveWorkbook

function New-DataSheet {
    param(
        [string]$SheetName,
        [object[]]$Rows
    )
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $SheetName
    for ($r = 0; $r -lt $Rows.Count; $r++) {
        $row = $Rows[$r]
        $newSheet.Cells.Item($r + 1, 1).Value = $row[0]
        $newSheet.Cells.Item($r + 1, 2).Value = $row[1]
        $newSheet.Cells.Item($r + 1, 3).Value = $row[2]
    }
}

$rows_gen6 = @(
    ,@(78735928, 718, 1025032)
    ,@(79347187, 726, 1011792)
    ,@(79058874, 723, 1021178)
    ,@(80404333, 737, 1004100)
    ,@(80500437, 738, 1000280)
    ,@(79058874, 723, 1024596)
    ,@(79539395, 728, 1016192)
    ,@(80212125, 735, 1015680)
    ,@(80404333, 737, 1005114)
    ,@(80884854, 742, 1000596)
)
New-DataSheet "gen_6" $rows_gen6

$rows_gen7 = @(
    ,@(78753245, 719, 1012814)
    ,@(78735928, 718, 1024760)
    ,@(79347187, 726, 1011472)
    ,@(80404333, 737, 1002548)
    ,@(80500437, 738, 1000280)
    ,@(79539395, 728, 1015984)
    ,@(79635499, 729, 1017564)
    ,@(80212125, 735, 1015316)
    ,@(79058874, 723, 1021164)
    ,@(80404333, 737, 1004100)
    ,@(80884854, 742, 1000596)
)
New-DataSheet "gen_7" $rows_gen7

$rows_gen8 = @(
    ,@(78753245, 719, 1012814)
    ,@(78735928, 718, 1024776)
    ,@(79347187, 726, 1011388)
    ,@(80404333, 737, 1002548)
    ,@(80500437, 738, 1000280)
    ,@(78753245, 719, 1016932)
    ,@(79347187, 726, 1011472)
    ,@(80884854, 742, 1000596)
    ,@(80404333, 737, 1003960)
    ,@(79539395, 728, 1015984)
    ,@(80212125, 735, 1015316)
    ,@(78753245, 719, 1024784)
    ,@(80404333, 737, 1004100)
)
New-DataSheet "gen_8" $rows_gen8

$rows_gen9 = @(
    ,@(78753245, 719, 1012814)
    ,@(80212125, 735, 999094)
    ,@(78735928, 718, 1024576)
    ,@(79347187, 726, 1011388)
    ,@(78753245, 719, 1016932)
    ,@(79347187, 726, 1011472)
    ,@(80404333, 737, 1002346)
    ,@(80500437, 738, 1000280)
    ,@(80884854, 742, 999792)
    ,@(80596542, 739, 1001008)
    ,@(79347187, 726, 1018020)
    ,@(80212125, 735, 1015316)
    ,@(80404333, 737, 1002548)
    ,@(80308229, 736, 1011892)
    ,@(79539395, 728, 1015984)
    ,@(78753245, 719, 1023860)
)
New-DataSheet "gen_9" $rows_gen9

$rows_gen10 = @(
    ,@(78753245, 719, 1012814)
    ,@(79347187, 726, 1004568)
    ,@(78735928, 718, 1025312)
    ,@(80212125, 735, 999094)
    ,@(78753245, 719, 1014668)
    ,@(79347187, 726, 1010156)
    ,@(80404333, 737, 1002346)
    ,@(80500437, 738, 1000280)
    ,@(80884854, 742, 999792)
    ,@(80404333, 737, 1002548)
    ,@(79347187, 726, 1011388)
    ,@(78753245, 719, 1016932)
)
New-DataSheet "gen_10" $rows_gen10

$wb.Worksheets.Item("gen_1").Activate()
Write-Output "Added sheets gen_6..gen_10"